$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6881498098373413
$ws.Range("B1").Value = 0.730939507484436
$ws.Range("C1").Value = 1.243634343147278
$ws.Range("D1").Value = 2.134048461914062
$ws.Range("E1").Value = 2.71191668510437
